$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("G3").Value = 1.95
$ws.Range("I3").Value = 4.5
$ws.Range("L3").Value = 5
$ws.Range("M3").Value = 1.13
$ws.Range("N3").Value = 6
$ws.Range("AC3").Value = 5.5
$ws.Range("AD3").Value = 8
$ws.Range("AF3").Value = 17
$ws.Range("AN3").Value = 9

# Row 4 updates
$ws.Range("G4").Value = 1.62
$ws.Range("I4").Value = 5
$ws.Range("S4").Value = 1.95
$ws.Range("T4").Value = 1.9
$ws.Range("W4").Value = 3.4
$ws.Range("X4").Value = 1.3
$ws.Range("AC4").Value = 7
$ws.Range("AH4").Value = 26
$ws.Range("AM4").Value = 301

# Row 5 updates
$ws.Range("G5").Value = 4.25
$ws.Range("H5").Value = 3.15
$ws.Range("I5").Value = 1.85
$ws.Range("J5").Value = 4.55
$ws.Range("L5").Value = 2.52
$ws.Range("P5").Value = 3
$ws.Range("S5").Value = 1.85
$ws.Range("X5").Value = 1.3
$ws.Range("Y5").Value = 1.42
$ws.Range("Z5").Value = 2.45
$ws.Range("AC5").Value = 12.5
$ws.Range("AI5").Value = 9.25
$ws.Range("AJ5").Value = 6.2
$ws.Range("AN5").Value = 6.8
$ws.Range("AO5").Value = 8.75
$ws.Range("AS5").Value = 26
